$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly added A/B values for existing rows 23 and 24
$ws.Range("A23").Value = 0.0312288
$ws.Range("B23").Value = 0.0293721

$ws.Range("A24").Value = 0.0213415
$ws.Range("B24").Value = 0.021437

# Add new rows 27-29 with C/D values (extends used range / dimension)
$ws.Range("C27").Value = 0.0072541
$ws.Range("D27").Value = 0.0037323

$ws.Range("C28").Value = 0.0201335
$ws.Range("D28").Value = 0.0077251

$ws.Range("C29").Value = 0.0092403
$ws.Range("D29").Value = 0.0047571
